$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.780.73"
$ws.Range("E2").Value = "  +3.59%  "

$ws.Range("D3").Value = "3.232.59"
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "541.61"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("D6").Value = "147.18"
$ws.Range("E6").Value = "  +5.22%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("E10").Value = "  +2.30%  "

$ws.Range("D11").Value = "0.438"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").Value = "3.788.41"
$ws.Range("E12").Value = "  +2.13%  "

$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("D14").Value = "26.13"
$ws.Range("E14").Value = "  +1.53%  "

$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("D16").Value = "60.788.85"
$ws.Range("E16").Value = "  +3.53%  "

$ws.Range("D17").Value = "3.230.70"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").Value = "6.32"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").Value = "13.36"
$ws.Range("E19").Value = "  +2.98%  "

$ws.Range("D20").Value = "8.36"
$ws.Range("E20").Value = "  +3.16%  "

$ws.Range("D21").Value = "376.76"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "0.527"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").Value = "70.03"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("E25").Value = "  +2.30%  "

$ws.Range("D26").Value = "8.69"
$ws.Range("E26").Value = "  +4.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").Value = "0.0₃0914"
$ws.Range("E28").Value = "  +5.88%  "

$ws.Range("D29").Value = "22.58"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").Value = "6.22"
$ws.Range("E31").Value = "  +3.07%  "

$ws.Range("D32").Value = "5.37"
$ws.Range("E32").Value = "  +4.35%  "

$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  +6.57%  "

$ws.Range("D34").Value = "6.62"
$ws.Range("E34").Value = "  +4.94%  "

$ws.Range("D35").Value = "158.31"
$ws.Range("E35").Value = "  +0.96%  "

$ws.Range("D36").Value = "1.41"
$ws.Range("E36").Value = "  +5.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.50"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.81%  "

$ws.Range("D38").Value = "2.803.60"
$ws.Range("E38").Value = "  +4.08%  "

$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("D40").Value = "0.0313"
$ws.Range("E40").Value = "  +8.16%  "

$ws.Range("E41").Value = "  +1.44%  "

$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("D43").Value = "39.99"
$ws.Range("E43").Value = "  +2.23%  "

$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "0.105"
$ws.Range("E45").Value = "  +1.78%  "

$ws.Range("D46").Value = "3.273.28"
$ws.Range("E46").Value = "  +2.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.26%  "

$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").Value = "20.97"
$ws.Range("E49").Value = "  +4.60%  "

$ws.Range("D50").Value = "0.802"
$ws.Range("E50").Value = "  +6.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "

